$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.263378620147705
$ws.Range("B1").Value = 2.403374433517456
$ws.Range("C1").Value = 4.323180675506592
$ws.Range("D1").Value = 2.636476516723633
$ws.Range("E1").Value = 1.347945213317871
